$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

# Copy formatting from row 14 (the row right below, which has correct original style) to new row 13
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(13).PasteSpecial(-4122)  # xlPasteFormats = -4122
$excel.CutCopyMode = $false

$dump = $ws.Range("C13").Worksheet.Name
Write-Host "done"
